{"js": "// The document holds a single table of \"two-digit \u00f7 one-digit\" drill\n// problems. Every 4th row (0, 4, 8, 12, 16 - 0 based) carries the five\n// exercise cells for that block; the three rows after each are blank\n// spacer rows. The edit swaps each exercise's expression text for a new\n// one while leaving the table shape (row/column counts) and per-cell\n// formatting untouched.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Row-major replacements: table row index -> new five cell values.\nconst replacements = {\n  0: [\"96\u00f74=\", \"49\u00f79=\", \"38\u00f75=\", \"67\u00f72=\", \"82\u00f72=\"],\n  4: [\"15\u00f79=\", \"78\u00f77=\", \"13\u00f72=\", \"63\u00f78=\", \"46\u00f76=\"],\n  8: [\"83\u00f77=\", \"17\u00f74=\", \"65\u00f74=\", \"44\u00f79=\", \"51\u00f72=\"],\n  12: [\"65\u00f73=\", \"32\u00f73=\", \"65\u00f72=\", \"36\u00f73=\", \"59\u00f78=\"],\n  16: [\"46\u00f75=\", \"80\u00f78=\", \"21\u00f76=\", \"86\u00f79=\", \"22\u00f73=\"],\n};\n\n// table.values expects the full grid; start from what's already there\n// and only overwrite the rows we care about so blank spacer rows (and\n// any other content) are left exactly as-is.\nconst newValues = table.values.map((row) => row.slice());\nfor (const key of Object.keys(replacements)) {\n  const idx = Number(key);\n  newValues[idx] = replacements[key];\n}\n\ntable.values = newValues;\nawait context.sync();\n", "ps1": "# The document holds a single table of \"two-digit \u00f7 one-digit\" drill\n# problems. Every 4th row (rows 1, 5, 9, 13, 17 in Word's 1-based Rows\n# collection) carries the five exercise cells for that block; the three\n# rows after each are blank spacer rows. The edit swaps each exercise's\n# expression text for a new one while leaving the table shape (row/column\n# counts) and per-cell formatting untouched.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Row-major replacements: 1-based table row -> new five cell values.\n$replacements = @{\n    1  = @(\"96\u00f74=\", \"49\u00f79=\", \"38\u00f75=\", \"67\u00f72=\", \"82\u00f72=\")\n    5  = @(\"15\u00f79=\", \"78\u00f77=\", \"13\u00f72=\", \"63\u00f78=\", \"46\u00f76=\")\n    9  = @(\"83\u00f77=\", \"17\u00f74=\", \"65\u00f74=\", \"44\u00f79=\", \"51\u00f72=\")\n    13 = @(\"65\u00f73=\", \"32\u00f73=\", \"65\u00f72=\", \"36\u00f73=\", \"59\u00f78=\")\n    17 = @(\"46\u00f75=\", \"80\u00f78=\", \"21\u00f76=\", \"86\u00f79=\", \"22\u00f73=\")\n}\n\nforeach ($rowIndex in $replacements.Keys) {\n    $values = $replacements[$rowIndex]\n    for ($col = 1; $col -le $values.Count; $col++) {\n        $t.Cell($rowIndex, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
